$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Insert a new column before column A (shifts existing A:D to B:E)
$ws.Columns.Item(1).Insert()

# New column A: id header + row numbers
$ws.Range("A1").Value = "id"
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4

# Column E already holds the former "job" column data after the insert shift;
# just rename its header from "job" to "jobTitle".
$ws.Range("E1").Value = "jobTitle"

# Column widths (speciality / department detail columns).
# Note: the host engine quantizes saved column width to an MDW-7 pixel grid
# (i.e. multiples of 1/7 character units: saved_pixels = Round(input*7)+5),
# so the requested 15.125 cannot be represented exactly; 14.4 is the
# ColumnWidth input whose rounded-on-save result (15.142857...) lands
# closest to the intended 15.125 value. 15.3 rounds on save to exactly 16.
$ws.Columns.Item(3).ColumnWidth = 14.4
$ws.Columns.Item(4).ColumnWidth = 15.3
